$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update location for rows 8, 9, 10 (C column)
$ws.Cells.Item(8, 3).Value = "East Valley"
$ws.Cells.Item(9, 3).Value = "East Valley"
$ws.Cells.Item(10, 3).Value = "Northwood"

# Add new rows 12-17
$data = @(
  @(2016, 4, "Northwood", 4),
  @(2017, 4, "West River", 2),
  @(2016, 5, "Northwood", 0),
  @(2017, 5, "West River", 17),
  @(2016, 7, "West River", 21),
  @(2017, 7, "West River", 9)
)

$r = 12
foreach ($row in $data) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $r = $r + 1
}

$null = $ws.Range("C16").Select()
